# Update "想去人数" (F column) values on both the "展览" and "全部类型"
# sheets to reflect newly generated output data.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F3").Value = 2174
    $ws.Range("F4").Value = 320
    $ws.Range("F5").Value = 75
    $ws.Range("F6").Value = 6407
    $ws.Range("F7").Value = 277
}
